$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sources")

# Copy the formatting of the last existing data row down onto the new row
# so the new row picks up the same cell styles (Hyperlink style for column A,
# wrap-text style for the rest) without minting any new style entries.
$ws.Range("A29:D29").Copy()
$ws.Range("A30:D30").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new source row
$ws.Range("A30").Value = "https://marinedebris.engr.uga.edu"
$ws.Range("B30").Value = "Debris traking and identification."
$ws.Range("C30").Value = "I will download the collection of plastic data and do some EDA in R."
$ws.Range("D30").Value = "Stuart"

# Grow the worksheet table (ListObject) so the new row is part of Table1
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:D30"))

# Update the view to where the new row was edited
$ws.Range("F30:F31").Select()

$wb.Save()
